# Generate Report for Handback
# Adds two new handed-back files (2f9fbfb4-6014-4175-bde6-ae24a26cf774 and
# 5a3ed23a-ad82-45b4-a5d5-a5d3b9d599ca) as new rows on the Overview, zh-cn
# and de-de worksheets, mirroring the layout of the existing rows.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$includeText = "Include"
$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)

$overviewRows = @(
    @{ Row = 6; File = "2f9fbfb4-6014-4175-bde6-ae24a26cf774.md" },
    @{ Row = 7; File = "5a3ed23a-ad82-45b4-a5d5-a5d3b9d599ca.md" }
)

foreach ($r in $overviewRows) {
    $row = $r.Row
    $file = $r.File

    $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$([System.IO.Path]::GetFileNameWithoutExtension($file))/e2e/$file"

    $wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($row, 1), $mdUrl, "", "", $file) | Out-Null
    $wsOverview.Cells.Item($row, 2).Value = $statusText
    $wsOverview.Cells.Item($row, 3).Value = $statusText
}

# ---------------------------------------------------------------------------
# Helper data describing the two language-specific detail sheets
# ---------------------------------------------------------------------------
$langSheets = @(
    @{
        SheetIndex = 2
        Lang       = "zh-cn"
        Rows       = @(
            @{
                Row            = 6
                MdFile         = "2f9fbfb4-6014-4175-bde6-ae24a26cf774.md"
                XlfFile        = "2f9fbfb4-6014-4175-bde6-ae24a26cf774.e68cf24d2f87f67fdc14a38318852c2a6d85c13b.zh-cn.xlf"
                HandoffTime    = "2016-01-27 08:48:13"
                HandbackTime   = "2016-01-27 08:48:58"
            },
            @{
                Row            = 7
                MdFile         = "5a3ed23a-ad82-45b4-a5d5-a5d3b9d599ca.md"
                XlfFile        = "5a3ed23a-ad82-45b4-a5d5-a5d3b9d599ca.5ad9f4e741d74b2569333a8b30945b7ba4827a5d.zh-cn.xlf"
                HandoffTime    = "2016-01-27 08:48:13"
                HandbackTime   = "2016-01-27 08:48:58"
            }
        )
    },
    @{
        SheetIndex = 3
        Lang       = "de-de"
        Rows       = @(
            @{
                Row            = 6
                MdFile         = "2f9fbfb4-6014-4175-bde6-ae24a26cf774.md"
                XlfFile        = "2f9fbfb4-6014-4175-bde6-ae24a26cf774.e68cf24d2f87f67fdc14a38318852c2a6d85c13b.de-de.xlf"
                HandoffTime    = "2016-01-27 08:48:25"
                HandbackTime   = "2016-01-27 08:49:20"
            },
            @{
                Row            = 7
                MdFile         = "5a3ed23a-ad82-45b4-a5d5-a5d3b9d599ca.md"
                XlfFile        = "5a3ed23a-ad82-45b4-a5d5-a5d3b9d599ca.5ad9f4e741d74b2569333a8b30945b7ba4827a5d.de-de.xlf"
                HandoffTime    = "2016-01-27 08:48:25"
                HandbackTime   = "2016-01-27 08:49:20"
            }
        )
    }
)

foreach ($langSheet in $langSheets) {
    $ws = $wb.Worksheets.Item($langSheet.SheetIndex)
    $lang = $langSheet.Lang

    foreach ($r in $langSheet.Rows) {
        $row = $r.Row
        $mdFile = $r.MdFile
        $xlfFile = $r.XlfFile
        $mdSha = [System.IO.Path]::GetFileNameWithoutExtension($mdFile)
        $xlfHash = $xlfFile.Split('.')[1]

        $mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/$mdSha/e2e/$mdFile"
        $handoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$xlfHash/ol-handoff/OpenLocalizationTestOrg/oltest.$lang/xinjiang/$xlfFile"
        $handbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/$xlfHash/ol-handback/OpenLocalizationTestOrg/oltest.$lang/xinjiang/$xlfFile"
        $targetUrl = "https://github.com/OpenLocalizationTestOrg/oltest.$lang/blob/$xlfHash/e2e/$mdFile"

        # Column A: Source File Name (hyperlink to the .md handoff source)
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 1), $mdUrl, "", "", $mdFile) | Out-Null

        # Column B: Status
        $ws.Cells.Item($row, 2).Value = $statusText

        # Column C: Correspond Handoff File (hyperlink to the handoff .xlf)
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 3), $handoffUrl, "", "", $xlfFile) | Out-Null

        # Column D: Correspond Handoff Datetime
        $ws.Cells.Item($row, 4).Value = $r.HandoffTime
        $ws.Cells.Item($row, 4).NumberFormat = $dateFmt

        # Column E: Target File (hyperlink to the .md file in the target-language repo)
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), $mdUrl, "", "", $mdFile) | Out-Null

        # Column F: Correspond Handback File (hyperlink to the handback .xlf)
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), $handbackUrl, "", "", $xlfFile) | Out-Null

        # Column G: Correspond Handback DateTime (plain text, matches existing rows' styling)
        $ws.Cells.Item($row, 7).Value = $r.HandbackTime

        # Column H: Handoff Reason
        $ws.Cells.Item($row, 8).Value = $includeText
    }
}
